$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'Given the predicted need for continued SARS-CoV-2 diagnostic testing, as well as the evolving availability and types of diagnostic tests, off-site COVID-19 testing centers (OSCTC) leaders need timely guidance to ensure they are meeting the needs of their unique populations.
 This research discusses the challenges and offers considerations for healthcare organizations and others when setting up and running OSCTCs.
 It also provides a springboard to engage policy makers and leaders in the healthcare community in a discussion about emergency preparedness, and how to better respond to testing needs going forward.
'
$ws.Range("E2").Value = '[Craig%Brammer%NULL%0, Stacy%Donohue%NULL%1, Timothy%Elwell%NULL%1, Eliza%Fishbein%NULL%1, D''Ante%Forschino%NULL%1, Dorothy%Horne%NULL%1, Buffy%Lloyd-Krejci%NULL%1, Jessica%Little%NULL%1, Bistra%Nikiforova%NULL%1, Elizabeth%Winterbauer%NULL%1]'
$ws.Range("I2").Value = ''
$ws.Range("J2").Value = 'The Authors. Published by Elsevier Inc.'

$ws.Range("D3").Value = 'Background
id="Par1">Early in the pandemic, inadequate SARS-CoV-2 testing limited understanding of transmission.

 Chief among barriers to large-scale testing was unknown feasibility, particularly in non-urban areas.

 Our objective was to report methods of high-volume, comprehensive SARS-CoV-2 testing, offering one model to augment disease surveillance in a rural community.


Methods
id="Par2">A community-university partnership created an operational site used to test most residents of Bolinas, California regardless of symptoms in 4 days (April 20th – April 23rd, 2020).

 Prior to testing, key preparatory elements included community mobilization, pre-registration, volunteer recruitment, and data management.

 On day of testing, participants were directed to a testing lane after site entry.

 An administrator viewed the lane-specific queue and pre-prepared test kits, linked to participants’ records.

 Medical personnel performed sample collection, which included finger prick with blood collection to run laboratory-based antibody testing and respiratory specimen collection for polymerase chain reaction (PCR).


Results
id="Par3">Using this 4-lane model, 1,840 participants were tested in 4 days.

 A median of 57 participants (IQR 47–67) were tested hourly.

 The fewest participants were tested on day 1 (n = 338 participants), an intentionally lower volume day, increasing to n = 571 participants on day 4. The number of testing teams was also increased to two per lane to allow simultaneous testing of multiple participants on days 2–4. Consistent staffing on all days helped optimize proficiency, and strong community partnership was essential from planning through execution.


Conclusions
id="Par4">High-volume ascertainment of SARS-CoV-2 prevalence by PCR and antibody testing was feasible when conducted in a community-led, drive-through model in a non-urban area.


Supplementary Information
The online version contains supplementary material available at 10.1186/s13690-021-00647-8.
'
$ws.Range("E3").Value = '[Ayesha%Appa%Ayesha.appa@ucsf.edu%0, Gabriel%Chamie%NULL%0, Aenor%Sawyer%NULL%1, Kimberly%Baltzell%NULL%1, Kathryn%Dippell%NULL%1, Salu%Ribeiro%NULL%1, Elias%Duarte%NULL%1, Joanna%Vinden%NULL%1, CLIAHUB%Consortium%NULL%1, Jonathan%Kramer-Feldman%NULL%1, Shahryar%Rahdari%NULL%1, Doug%MacIntosh%NULL%1, Katherine%Nicholson%NULL%1, Jonathan%Im%NULL%1, Diane%Havlir%NULL%1, Bryan%Greenhouse%NULL%0]'
$ws.Range("I3").Value = ''
$ws.Range("J3").Value = 'BioMed Central'

$ws.Range("D4").Value = 'In this study, we evaluated the efficiency of a drive-through (DT) screening system for severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) by comparing it with a conventional screening system.
 We reviewed and analyzed the SARS-CoV-2 screening data obtained at our university hospital.
 We compared the number of tests for SARS-CoV-2 (using real-time polymerase chain reaction) performed using two different specimen collection systems—DT and conventional—during the coronavirus disease 2019 (COVID-19) outbreak in Daegu.
 Based on the results, the DT screening system collected 5.8 times more specimens for testing than the conventional screening system.
 From 27 January to 31 March 2020, 6211 individuals were screened for SARS-CoV-2 infection using either the DT or conventional system.
 In total, 217 individuals tested positive for SARS-CoV-2 (positive rate: 3.50%).
 Of the 6211 individuals, 3368 were symptomatic or had a history of contact with COVID-19 patients, and 142 of them tested positive for SARS-CoV-2 (positive rate: 4.22%).
 Further, 2843 individuals were asymptomatic and had no history of contact with COVID-19 patients, and 75 of them tested positive for SARS-CoV-2 (positive rate: 2.64%).
 In conclusion, the DT system allowed clinicians to collect specimens for SARS-CoV-2 screening more efficiently than the conventional system.
 Furthermore, as there might be several COVID-19 patients who remain asymptomatic, expanding the screening test to asymptomatic individuals would be necessary.
'
$ws.Range("E4").Value = '[Min Cheol%Chang%NULL%0, Wan-Seok%Seo%NULL%1, Donghwi%Park%NULL%0, Jian%Hur%NULL%0]'
$ws.Range("I4").Value = ''
$ws.Range("J4").Value = 'MDPI'

$ws.Range("E5").Value = '[Marci L.%Drees%NULL%0, Mia A.%Papas%NULL%1, Terri E.%Corbo%NULL%1, Kimberly D.%Williams%NULL%1, Sharon T.%Kurfuerst%NULL%2, Sharon T.%Kurfuerst%NULL%0]'
$ws.Range("I5").Value = ''
$ws.Range("J5").Value = 'Cambridge University Press'

$ws.Range("E6").Value = '[Erin F.%Flynn%flynne@email.chop.edu%0, Elizabeth%Kuhn%NULL%1, Mohammed%Shaik%NULL%1, Elizabeth%Tarr%NULL%1, Nicole%Scattolini%NULL%1, Allison%Ballantine%NULL%1]'
$ws.Range("I6").Value = ''
$ws.Range("J6").Value = 'by Academic Pediatric Association'

$ws.Range("E7").Value = '[Scott A.%Goldberg%NULL%0, Robert A.%Bonacci%NULL%1, Lucas C.%Carlson%NULL%1, Charles T.%Pu%NULL%1, Christine S.%Ritchie%NULL%1]'
$ws.Range("I7").Value = ''
$ws.Range("J7").Value = 'Department of Emergency Medicine, University of California, Irvine School of Medicine'

$ws.Range("E8").Value = '[Travis%Sanchez%NULL%0, Sadhu%Panda%NULL%2, Sadhu%Panda%NULL%0, Ebrahim%Khajeh%NULL%1, Alexandra%Halalau%alexandra.halalau@beaumont.edu%2, Alexandra%Halalau%alexandra.halalau@beaumont.edu%0, Jeffrey%Ditkoff%NULL%2, Jeffrey%Ditkoff%NULL%0, Jessica%Hamilton%NULL%2, Jessica%Hamilton%NULL%0, Aryana%Sharrak%NULL%2, Aryana%Sharrak%NULL%0, Aimen%Vanood%NULL%2, Aimen%Vanood%NULL%0, Amr%Abbas%NULL%2, Amr%Abbas%NULL%0, James%Ziadeh%NULL%2, James%Ziadeh%NULL%0]'
$ws.Range("I8").Value = ''
$ws.Range("J8").Value = 'JMIR Publications'

$ws.Range("D9").Value = 'To increase the country’s capacity to test and track suspected coronavirus disease 2019 (COVID-19) cases, Israel launched drive-through testing centers in key cities, including Tel Aviv, Jerusalem, Be’er Sheva, and Haifa.
 This article examines the challenges that the national emergency medical services and volunteers faced in the process of implementing drive-through testing centers to offer lessons learned and direction to health-care professionals in other countries.
'
$ws.Range("E9").Value = '[Edward%Kim%NULL%0]'
$ws.Range("I9").Value = ''
$ws.Range("J9").Value = 'Cambridge University Press'

$ws.Range("D10").Value = 'With the ongoing novel coronavirus disease 2019 (COVID-19) pandemic, the number of individuals that need to be tested for COVID-19 has been rapidly increasing.
 A walk-through (WT) screening center using negative pressure booths that is inspired by the biosafety cabinet has been designed and implemented in Korea for easy screening of COVID-19 and for safe and efficient consultation for patients with fever or respiratory symptoms.
 Here, we present the overall concept, advantages, and limitations of the COVID-19 WT screening center.
 The WT center increases patient access to the screening clinics and adequately protects healthcare personnel while reducing the consumption of personal protective equipment.
 It can also increase the number of people tested by 9–10 fold.
 However, there is a risk of cross-infection at each stage of screening treatment, including the booths, and adverse reactions with disinfection of the booths.
 These limitations can be overcome using mobile technology and increasing the number of booths to reduce congestion inside the center, reducing booth volume for sufficient and rapid ventilation, and using an effective, harmless, and certified environmental disinfectant.
 A WT center can be implemented in other institutions and countries and modified depending on local needs to cope with the COVID-19 pandemic.
'
$ws.Range("E10").Value = '[Sang Il%Kim%NULL%0, Ji Yong%Lee%NULL%2, Ji Yong%Lee%NULL%0]'
$ws.Range("I10").Value = ''
$ws.Range("J10").Value = 'The Korean Academy of Medical Sciences'

$ws.Range("C11").Value = 'Unknown Title'
$ws.Range("E11").Value = '[]'
$ws.Range("F11").Value = 'not found'
$ws.Range("G11").Value = 'N/A'
$ws.Range("H11").Value = '1970-01-01'
$ws.Range("J11").Value = ''

$ws.Range("D12").Value = 'As the coronavirus disease 2019 (COVID-19) outbreak is ongoing, the number of individuals to be tested for COVID-19 is rapidly increasing.
 For safe and efficient screening for COVID-19, drive-through (DT) screening centers have been designed and implemented in Korea.
 Herein, we present the overall concept, advantages, and limitations of the COVID-19 DT screening centers.
 The steps of the DT centers include registration, examination, specimen collection, and instructions.
 The entire service takes about 10 minutes for one testee without leaving his or her cars.
 Increased testing capacity over 100 tests per day and prevention of cross-infection between testees in the waiting space are the major advantages, while protection of staff from the outdoor atmosphere is challenging.
 It could be implemented in other countries to cope with the global COVID-19 outbreak and transformed according to their own situations.
'
$ws.Range("E12").Value = '[Ki Tae%Kwon%NULL%0, Jae-Hoon%Ko%NULL%2, Jae-Hoon%Ko%NULL%0, Heejun%Shin%NULL%2, Heejun%Shin%NULL%0, Minki%Sung%NULL%2, Minki%Sung%NULL%0, Jin Yong%Kim%NULL%0, Jin Yong%Kim%NULL%0]'
$ws.Range("I12").Value = ''
$ws.Range("J12").Value = 'The Korean Academy of Medical Sciences'

$ws.Range("D13").Value = 'As the world witnessed the rapid spread of SARS-CoV-2, the World Health Organization has called for governing bodies worldwide to intensify case findings, contact tracing, monitoring, and quarantine or isolation of contacts with COVID-19. Drive-through (DT) screening is a form of case detection which has recently gain preference globally.
 Proper implementation of this system can help remediate the outbreak.
'
$ws.Range("E13").Value = '[Elly%Lee%NULL%0, Nurul Yaqeen%Mohd Esa%NULL%1, Tong Ming%Wee%NULL%1, Chun Ian%Soo%NULL%1]'
$ws.Range("I13").Value = ''
$ws.Range("J13").Value = 'Taiwan Society of Microbiology. Published by Elsevier Taiwan LLC.'

$ws.Range("D14").Value = 'There is still a paucity of studies on real-world outcome of screening clinic for hospital protection from coronavirus disease 2019 (COVID-19).
 As the number of COVID-19 cases was growing rapidly in Daegu, Korea, we started operating an active screening clinic outside the hospital premises.
 Over two weeks, 2,087 patients were screened using real-time reverse transcriptase polymerase chain reaction testing for severe acute respiratory syndrome coronavirus 2, with 42 confirmed cases.
 Before the screening clinic period, an average of 36 beds (maximum 67 beds) per day were closed due to unrecognized COVID-19 patients entering the hospital.
 In contrast, after the screening clinic operated well, only one event of closing emergency room (25 beds) occurred due to a confirmed COVID-19 case of asymptomatic patient.
 We report the operational process of screening clinic for COVID-19 and its effectiveness in maintaining the function of tertiary hospitals.
'
$ws.Range("E14").Value = '[Yong Shik%Kwon%NULL%0, Sun Hyo%Park%NULL%2, Sun Hyo%Park%NULL%0, Hyun Jung%Kim%NULL%2, Hyun Jung%Kim%NULL%0, Ji Yeon%Lee%NULL%0, Ji Yeon%Lee%NULL%0, Mi-ri%Hyun%NULL%2, Mi-ri%Hyun%NULL%0, Hyun ah%Kim%NULL%2, Hyun ah%Kim%NULL%0, Jae Seok%Park%NULL%0, Jae Seok%Park%NULL%0]'
$ws.Range("I14").Value = ''
$ws.Range("J14").Value = 'The Korean Academy of Medical Sciences'

$ws.Range("D15").Value = 'id="Par1">In Taiwan, high-risk patients have been identified and tested for preventing community spread of COVID-19. Most sample collection was performed in emergency departments (EDs).
 Traditional sample collection requires substantial personal protective equipment (PPE), healthcare professionals, sanitation workers, and isolation space.
 To solve this problem, we established a multifunctional sample collection station (MSCS) for COVID-19 testing in front of our ED.
 The station is composed of a thick and clear acrylic board (2 cm), which completely separates the patient and medical personnel.
 Three pairs of gloves (length, 45 cm) are attached and fixed on the outside wall of the MSCS.
 The gloves are used to conduct sampling of throat/nasal swabs, sputum, and blood from patients.
 The gap between the board and the building is only 0.2 cm (sealed with silicone sealant).
 ED personnel communicate with patients using a small two-way broadcast system.
 Medical waste is put in specific trashcans installed in the table outside the MSCS.
 With full physical protection, the personnel conducting the sampling procedure need to wear only their N95 mask and gloves.
 After we activated the station, our PPE, sampling time, and sanitization resources were considerably conserved during the 4-week observation period.
 The MSCS obviously saved time and PPE.
 It elevated the efficiency and capacity of the ED for handling potential community infections of COVID-19.'
$ws.Range("E15").Value = '[Po-Ting%Lin%NULL%0, Ting-Yuan%Ni%NULL%1, Tren-Yi%Chen%NULL%1, Chih-Pei%Su%NULL%1, Hsiao-Fen%Sun%NULL%1, Mu-Kuan%Chen%NULL%1, Chu-Chung%Chou%NULL%1, Po-Yu%Wang%NULL%1, Yan-Ren%Lin%H6213.lac@gmail.com%1]'
$ws.Range("I15").Value = ''
$ws.Range("J15").Value = 'BioMed Central'

$ws.Range("D16").Value = 'Drive-through coronavirus disease 2019 screening can evaluate large numbers of patients while reducing healthcare exposures and personal protective equipment use.
 We describe the characteristics of screened individuals as well as drive-through process and outcome measures.
 Optimal drive-through screening involves rapid turnaround of test results and linkage to follow-up care.
'
$ws.Range("E16").Value = '[David A%Lindholm%david.a.lindholm4.mil@mail.mil%0, John L%Kiley%NULL%2, John L%Kiley%NULL%0, Nathan K%Jansen%NULL%1, Robert T%Hoard%NULL%1, Matthew R%Bondaryk%NULL%1, Elizabeth M%Stanley%NULL%1, Gadiel R%Alvarado%NULL%1, Ana E%Markelz%NULL%1, Robert J%Cybulski%NULL%1, Jason F%Okulicz%NULL%1]'
$ws.Range("I16").Value = ''
$ws.Range("J16").Value = 'Oxford University Press'

$ws.Range("D17").Value = 'To combat the ongoing COVID‐19 pandemic, Singapore has adopted a rigorous screening approach that involves aggressive contact tracing, rapid isolation of confirmed or suspect cases, and immediate ring‐fencing of emerging local clusters and hotspots.
 Our screening centre team has been involved in running Singapore''s designated screening centre since the end of January this year.
 With a well‐defined blueprint and substantial pre‐outbreak preparatory work, initial operations at our screening centre commenced within a day on activation and full operational status was attained in 3 days.
 As of 8 April 2020, the screening centre had screened more than 14,000 patients.
 We have adopted a “whole‐of‐hospital” approach, enlisting the help from other departments and subspecialties to augment manpower.
 Meticulous infrastructure planning to facilitate patient flow and strict measures to prevent nosocomial transmission and occupational exposure were instituted to safeguard both the staff and patients.
 This paper aims to describe our key takeaways in the course of operations and discuss the challenges encountered.
'
$ws.Range("E17").Value = '[Charmaine Malenab%Manauis%NULL%0, Marvin%Loh%marvin.loh@mohh.com.sg%1, James%Kwan%NULL%1, John%Chua Mingzhou%NULL%1, Han Jie%Teo%NULL%1, David%Teng Kuan Peng%NULL%1, Shawn%Vasoo Sushilan%NULL%1, Yee Sin%Leo%NULL%1, Ang%Hou%NULL%1]'
$ws.Range("I17").Value = ''
$ws.Range("J17").Value = 'John Wiley and Sons Inc.'

$ws.Range("C18").Value = 'Unknown Title'
$ws.Range("D18").Value = 'Unknown Abstract'
$ws.Range("E18").Value = '[]'
$ws.Range("F18").Value = 'not found'
$ws.Range("G18").Value = 'N/A'
$ws.Range("I18").Value = ''

$ws.Range("D19").Value = 'The ongoing coronavirus disease 2019 (COVID-19) pandemic is causing tremendous damage globally.
 The Republic of Korea (ROK), a highly export-dependent nation, is a leader in the fight against the COVID-19 pandemic and coping well with the disaster.
 Like the drive-through COVID-19 testing, which reflects the brilliant flexibility of the Korean medical system, onsite mass workplace testing for COVID-19, which our hospital has been performing over the past few months, is a unique and valuable countermeasure.
 We believe it is time that the current health examination system for workers in the ROK considered the risk of transmissible diseases.
'
$ws.Range("E19").Value = '[Eunhye%Seo%NULL%0, Eunchan%Mun%NULL%2, Eunchan%Mun%NULL%0, Wonsool%Kim%NULL%2, Wonsool%Kim%NULL%0, Changhwan%Lee%NULL%2, Changhwan%Lee%NULL%0]'
$ws.Range("I19").Value = ''
$ws.Range("J19").Value = 'Korean Society of Occupational & Environmental Medicine'

$ws.Range("D20").Value = 'The World Health Organization declared COVID-19 a global pandemic in March 2020. A major challenge in this worldwide pandemic has been efficient and effective large-scale testing for the disease.
 In this communication, we discuss lessons learned in the set up and function of a locally organized drive-through testing facility.
'
$ws.Range("E20").Value = '[Aditya%Shah%shah.aditya@mayo.edu%0, Douglas%Challener%NULL%1, Aaron J.%Tande%NULL%1, Maryam%Mahmood%NULL%1, John C.%O’Horo%NULL%1, Elie%Berbari%NULL%1, Sarah J.%Crane%NULL%1]'
$ws.Range("I20").Value = ''
$ws.Range("J20").Value = 'Mayo Foundation for Medical Education and Research'

$ws.Range("E21").Value = '[Angie N.%Ton%NULL%0, Tarang%Jethwa%NULL%1, Karen%Waters%NULL%2, Karen%Waters%NULL%0, Leigh L.%Speicher%NULL%2, Leigh L.%Speicher%NULL%0, Dawn%Francis%NULL%2, Dawn%Francis%NULL%0]'
$ws.Range("I21").Value = ''
$ws.Range("J21").Value = 'Association for Professionals in Infection Control and Epidemiology, Inc. Published by Elsevier Inc.'

$ws.Range("D22").Value = 'In response to the outbreak of COVID-19, we set up a team to carry out sampling in the community.
 This enabled individuals to remain in self-isolation in their own homes and to prevent healthcare settings and services from being overwhelmed by admissions for sampling of suspected cases.
 There is evidence that this is a cost effective, safe and necessary service to complement COVID-19 testing in hospitals.
'
$ws.Range("E22").Value = '[Kate%Mark%NULL%0, Katie%Steel%NULL%1, Janet%Stevenson%NULL%1, Christine%Evans%NULL%1, Duncan%McCormick%NULL%1, Lorna%Willocks%NULL%1, Alison%McCallum%NULL%1, Laura%Jones%NULL%1, Ingolfur%Johannessen%NULL%1, Kate%Templeton%NULL%1, Oliver%Koch%NULL%0, Claire%Mackintosh%NULL%1]'
$ws.Range("I22").Value = ''
$ws.Range("J22").Value = 'European Centre for Disease Prevention and Control (ECDC)'

